$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that immediately follows the
# "Klienta: ... proźbą o serwis." bullet - that's the paragraph we need
# to turn into the new bold/italic/dark-red note about the .env file.
$i = 0
$targetIndex = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -match "prozba o serwis" -or $t -match "proźbą o serwis") {
        $targetIndex = $i + 1
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range

# Prime the (still empty) paragraph mark with the new formatting first -
# this is what makes the paragraph-mark run properties (w:pPr/w:rPr) pick
# up bCs/iCs alongside b/i/color once text is typed in below.
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Italic = $true
$r.Font.ItalicBi = $true
$r.Font.Color = 192

# Insert the new sentence (kept as four consecutive inserts, mirroring
# the four runs in the source edit - Word will coalesce runs that end up
# with identical formatting, exactly like the authoring app did).
$r.InsertAfter("Dodatkowym elementem jest ustawienie pliku konfiguracyjnego .env by móc resetować lub ")

$p1 = $d.Paragraphs.Item($targetIndex)
$ins = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$ins.InsertAfter("odzyskiwać ")

$p2 = $d.Paragraphs.Item($targetIndex)
$ins = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$ins.InsertAfter("hasło")

$p3 = $d.Paragraphs.Item($targetIndex)
$ins = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$ins.InsertAfter(". W pliku .env podajemy dane testowego klienta pocztowego by móc otrzymać hasło")

# Apply bold / italic / dark-red formatting to the whole paragraph
# (covers both the runs and the paragraph-mark run properties).
$final = $d.Paragraphs.Item($targetIndex)
$fr = $final.Range
$fr.Font.Bold = $true
$fr.Font.BoldBi = $true
$fr.Font.Italic = $true
$fr.Font.ItalicBi = $true
$fr.Font.Color = 192

Write-Host "Paragraph $targetIndex updated: [$($final.Range.Text)]"
